$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.199.23'
$ws.Range("E2").Value = '  -0.10%  '

$ws.Range("D3").Value = '1.850.72'
$ws.Range("E3").Value = '  -0.44%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.39'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6990'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.85%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07724'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.34%  '

$ws.Range("E9").Value = '  -1.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.55'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.17%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07824'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '93.21'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.24%  '

$ws.Range("D13").Value = '1.847.01'
$ws.Range("E13").Value = '  -0.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.133'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.76%  '

$ws.Range("E15").Value = '  -0.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.650'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.29%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008322'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.30%  '

$ws.Range("D18").Value = '29.197.32'
$ws.Range("E18").Value = '  -0.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '241.50'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.50%  '

$ws.Range("D20").Value = '2.087.08'
$ws.Range("E20").Value = '  -1.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.74'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.530'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1511'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.97'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.94%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.846'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.29'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.543'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.232'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.175'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.192'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05120'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.60%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.8005'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +5.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.871'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.49%  '

$ws.Range("E36").Value = '  -1.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.691'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.58%  '

$ws.Range("D38").Value = '1.312.73'
$ws.Range("E38").Value = '  +8.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01873'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.714'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9442'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +5.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.024'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +8.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '107.10'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.51%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.0000'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.730'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.43%  '

$ws.Range("E46").Value = '  +1.32%  '

$ws.Range("D47").Value = '1.988.60'
$ws.Range("E47").Value = '  -0.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5178'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '64.18'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.64%  '

$ws.Range("E50").Value = '  +0.98%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.999'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.11%  '

